$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.561.03"
$ws.Range("E2").Value = "'  +0.37%  "

$ws.Range("D3").Value = "'3.253.91"
$ws.Range("E3").Value = "'  -2.74%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.34%  "

$ws.Range("D5").Value = "'576.07"
$ws.Range("E5").Value = "'  -1.57%  "

$ws.Range("D6").Value = "'170.72"
$ws.Range("E6").Value = "'  -8.13%  "

$ws.Range("E7").Value = "'  +0.06%  "

$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "'  -0.38%  "

$ws.Range("D9").Value = "'3.248.26"
$ws.Range("E9").Value = "'  -2.74%  "

$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "'  -6.57%  "

$ws.Range("D11").Value = "'0.564"
$ws.Range("E11").Value = "'  -3.12%  "

$ws.Range("D12").Value = "'44.37"
$ws.Range("E12").Value = "'  -5.60%  "

$ws.Range("D13").Value = "'0.0000266"
$ws.Range("E13").Value = "'  -1.06%  "

$ws.Range("D14").Value = "'683.34"
$ws.Range("E14").Value = "'  +4.22%  "

$ws.Range("D15").Value = "'3.789.62"
$ws.Range("E15").Value = "'  +4.31%  "

$ws.Range("D16").Value = "'8.15"
$ws.Range("E16").Value = "'  -4.11%  "

$ws.Range("D17").Value = "'66.728.17"
$ws.Range("E17").Value = "'  +0.48%  "

$ws.Range("E18").Value = "'  +0.23%  "

$ws.Range("D19").Value = "'3.259.71"
$ws.Range("E19").Value = "'  -2.31%  "

$ws.Range("D20").Value = "'17.04"
$ws.Range("E20").Value = "'  -4.75%  "

$ws.Range("D21").Value = "'10.53"
$ws.Range("E21").Value = "'  -5.39%  "

$ws.Range("D22").Value = "'0.873"
$ws.Range("E22").Value = "'  -2.96%  "

$ws.Range("D23").Value = "'16.74"
$ws.Range("E23").Value = "'  -5.38%  "

$ws.Range("D24").Value = "'5.20"
$ws.Range("E24").Value = "'  +2.76%  "

$ws.Range("D25").Value = "'96.97"
$ws.Range("E25").Value = "'  -3.33%  "

$ws.Range("D26").Value = "'3.81"
$ws.Range("E26").Value = "'  -4.81%  "

$ws.Range("D27").Value = "'2.61"
$ws.Range("E27").Value = "'  -6.70%  "

$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "'  -7.68%  "

$ws.Range("D29").Value = "'32.18"
$ws.Range("E29").Value = "'  +0.27%  "

$ws.Range("D30").Value = "'8.20"
$ws.Range("E30").Value = "'  -4.11%  "

$ws.Range("D31").Value = "'6.59"
$ws.Range("E31").Value = "'  -3.73%  "

$ws.Range("D32").Value = "'566.62"
$ws.Range("E32").Value = "'  -5.99%  "

$ws.Range("D33").Value = "'10.74"
$ws.Range("E33").Value = "'  -3.20%  "

$ws.Range("D34").Value = "'3.773.10"
$ws.Range("E34").Value = "'  -2.75%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "'  -0.07%  "

$ws.Range("E36").Value = "'  -4.32%  "

$ws.Range("D37").Value = "'54.87"
$ws.Range("E37").Value = "'  -2.88%  "

$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "'  -16.85%  "

$ws.Range("B40").Value = "'Fetch.AI"
$ws.Range("C40").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'2.53"
$ws.Range("E40").Value = "'  -7.42%  "

$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'31.03"
$ws.Range("E41").Value = "'  -6.24%  "

$ws.Range("B42").Value = "'ApeXProtocol"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'3.25"
$ws.Range("E42").Value = "'  -3.48%  "

$ws.Range("D43").Value = "'0.0₃0647"
$ws.Range("E43").Value = "'  -7.81%  "

$ws.Range("D44").Value = "'0.320"
$ws.Range("E44").Value = "'  -6.22%  "

$ws.Range("D45").Value = "'2.93"
$ws.Range("E45").Value = "'  -8.30%  "

$ws.Range("D46").Value = "'0.0397"
$ws.Range("E46").Value = "'  -4.88%  "

$ws.Range("E47").Value = "'  +0.07%  "

$ws.Range("D48").Value = "'0.125"
$ws.Range("E48").Value = "'  -2.22%  "

$ws.Range("D49").Value = "'2.48"
$ws.Range("E49").Value = "'  -2.83%  "

$ws.Range("D50").Value = "'1.31"
$ws.Range("E50").Value = "'  -1.48%  "

$ws.Range("D51").Value = "'126.43"
$ws.Range("E51").Value = "'  -2.87%  "
